# Sprint 100 manual testcases update
# - Row 23 gets new Steps/Expected-output text (a "new" Options menu item test case)
#   and its row height grows to fit the longer wrapped text.
# - The active sheet view's selection moves to D23 (the edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 23 height -------------------------------------------------
$ws.Rows.Item(23).RowHeight = 95.25

# ---- D23 : "Steps" column -------------------------------------------
$d23 = $ws.Range("D23")
$d23.Value = ' Action->View suppliers->Marketlist screen->Options(new)'
# bold the "Action->View suppliers->Marketlist screen->Options" portion
$d23.Characters(2, 50).Font.Bold = $true

# ---- E23 : "Expected output" column ---------------------------------
$e23 = $ws.Range("E23")
$e23.Value = 'Select the ''Copy data to another outlet'',it comes pop up page is like                                                                                                                                  1. Select what to copy : Custom name, Tags, Buyer product code, Unit price and MOQ                                                                                                                    2.Copy to : Outlet name and Ok'
$e23.Characters(12, 29).Font.Bold = $true
$e23.Characters(200, 86).Font.Bold = $true
$e23.Characters(398, 30).Font.Bold = $true

# ---- Sheet view: move the selection to the edited cell --------------
$ws.Range("D23").Select()
